$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 0.2312925170068027
    "C2" = 0.5
    "J2" = 0.01360544217687075
    "P2" = 0.1496598639455782
    "S2" = 0.1054421768707483
    "B3" = 0.006211180124223602
    "C3" = 0.03726708074534162
    "J3" = 0.04968944099378882
    "P3" = 0.7453416149068323
    "S3" = 0.1614906832298137
    "J4" = 0.08695652173913043
    "P4" = 0.7608695652173914
    "S4" = 0.1521739130434783
    "B6" = 0.05045871559633028
    "D6" = 0.01834862385321101
    "F6" = 0.04128440366972477
    "J6" = 0.2752293577981652
    "O6" = 0.03211009174311927
    "Q6" = 0.1422018348623853
    "R6" = 0.0871559633027523
    "S6" = 0.3532110091743119
    "B7" = 0.1111111111111111
    "D7" = 0.0303030303030303
    "F7" = 0.04545454545454546
    "J7" = 0.1313131313131313
    "O7" = 0.005050505050505051
    "Q7" = 0.2424242424242424
    "R7" = 0.1161616161616162
    "S7" = 0.3181818181818182
    "B8" = 0.1228070175438596
    "D8" = 0.02255639097744361
    "F8" = 0.05012531328320802
    "J8" = 0.1027568922305764
    "O8" = 0.0100250626566416
    "Q8" = 0.2155388471177945
    "R8" = 0.09523809523809523
    "S8" = 0.3809523809523809
    "B9" = 0.1064814814814815
    "D9" = 0.009259259259259259
    "F9" = 0.06481481481481481
    "J9" = 0.1435185185185185
    "O9" = 0.01388888888888889
    "Q9" = 0.2175925925925926
    "R9" = 0.09259259259259259
    "S9" = 0.3518518518518519
    "B10" = 0.08719135802469136
    "D10" = 0.01929012345679012
    "F10" = 0.07561728395061729
    "J10" = 0.1419753086419753
    "O10" = 0.01157407407407407
    "Q10" = 0.2376543209876543
    "R10" = 0.08950617283950617
    "S10" = 0.3371913580246914
    "G11" = 0.1370716510903427
    "J11" = 0.1090342679127726
    "K11" = 0.1869158878504673
    "L11" = 0.5607476635514018
    "S11" = 0.006230529595015576
    "G12" = 0.7037037037037037
    "J12" = 0.2116402116402116
    "K12" = 0.02116402116402116
    "L12" = 0.04761904761904762
    "S12" = 0.01587301587301587
    "G13" = 0.8857142857142857
    "J13" = 0.1142857142857143
    "F15" = 0.03827751196172249
    "H15" = 0.1531100478468899
    "I15" = 0.09090909090909091
    "J15" = 0.3684210526315789
    "K15" = 0.0430622009569378
    "M15" = 0.004784688995215311
    "N15" = 0.009569377990430622
    "O15" = 0.05263157894736842
    "S15" = 0.2392344497607655
    "F16" = 0.0154639175257732
    "H16" = 0.1597938144329897
    "I16" = 0.1082474226804124
    "J16" = 0.3298969072164948
    "K16" = 0.1288659793814433
    "M16" = 0.01030927835051546
    "N16" = 0.005154639175257732
    "O16" = 0.07731958762886598
    "S16" = 0.1649484536082474
    "F17" = 0.01747572815533981
    "H17" = 0.1495145631067961
    "I17" = 0.0912621359223301
    "J17" = 0.4252427184466019
    "K17" = 0.1067961165048544
    "M17" = 0.01359223300970874
    "O17" = 0.05825242718446602
    "S17" = 0.1378640776699029
    "F18" = 0.01382488479262673
    "H18" = 0.1566820276497696
    "I18" = 0.06912442396313365
    "J18" = 0.4147465437788018
    "K18" = 0.1244239631336406
    "M18" = 0.03686635944700461
    "O18" = 0.06451612903225806
    "S18" = 0.119815668202765
    "F19" = 0.01751459549624687
    "H19" = 0.1843202668890742
    "I19" = 0.09257714762301918
    "J19" = 0.3653044203502919
    "K19" = 0.1142618849040867
    "M19" = 0.0158465387823186
    "N19" = 0.00250208507089241
    "O19" = 0.07172643869891576
    "S19" = 0.1359466221851543
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

Write-Output "Applied $($changes.Count) cell updates"
